# propertyEnumCast.xlsx — add the "theStepNam" row (name() accessor) below
# the existing cast/toString rows, and fix the AK,AS -> AK typo in F9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F9 previously held "AK,AS" — correct it to "AK".
$ws.Range("F9").Value = "AK"

# New row 13, mirroring the layout of rows 11/12 (label in D, formula-looking
# literal text in E, merged across E:F).
$ws.Range("D13").Style = "Normal"
$ws.Range("D13").Value = "theStepNam"
$ws.Range("E13").Value = "'=`$properties.state.name()"
$ws.Range("F13").Value = ""
$ws.Range("E13:F13").Merge()

Write-Host "done"
